$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values to match the refreshed cryptocurrency data feed.
# Cells whose new text looks like a plain number are forced to stay as
# text (matching the original inline-string cell type) by temporarily
# applying a text number format, then resetting the cell style so no
# extra formatting is left behind.

$ws.Range('D2').Value = '67.554.06'
$ws.Range('E2').Value = '  -1.20%  '
$ws.Range('D3').Value = '3.764.30'
$ws.Range('E3').Value = '  +0.55%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '594.39'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.12%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '166.86'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.32%  '
$ws.Range('D7').Value = '3.762.37'
$ws.Range('E7').Value = '  +0.53%  '
$ws.Range('E8').Value = '  +0.09%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.518'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.15%  '
$ws.Range('E10').Value = '  -0.21%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.29'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.28%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.446'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.04%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000253'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.27%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '35.99'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.06%  '
$ws.Range('D15').Value = '4.395.94'
$ws.Range('E15').Value = '  +0.58%  '
$ws.Range('D16').Value = '3.792.56'
$ws.Range('E16').Value = '  +1.85%  '
$ws.Range('D17').Value = '67.502.39'
$ws.Range('E17').Value = '  -1.19%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '18.33'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.43%  '
$ws.Range('E19').Value = '  +0.08%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.97'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.08%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.97'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -7.04%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '455.42'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.68%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.693'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.07%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.0000152'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +5.01%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '83.05'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.31%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.88'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.08%  '
$ws.Range('E27').Value = '  -1.88%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.10'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.03%  '
$ws.Range('E29').Value = '  +0.10%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.76'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.05%  '
$ws.Range('E31').Value = '  +2.71%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.23'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.33%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '29.54'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.95%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '9.13'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.07%  '
$ws.Range('E35').Value = '  +0.04%  '
$ws.Range('D36').Value = '3.717.57'
$ws.Range('E36').Value = '  +0.59%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0999'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.35%  '
$ws.Range('B38').Value = 'dogwifhat'
$ws.Range('C38').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.31'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.96%  '
$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.138'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.46%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.996'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.39%  '
$ws.Range('E41').Value = '  -0.60%  '
$ws.Range('E42').Value = '  -0.10%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '45.75'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +6.85%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '48.72'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +5.70%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.298'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.72%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.29'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.04%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '147.86'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.09%  '
$ws.Range('B49').Value = 'Bittensor'
$ws.Range('C49').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '388.03'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.15%  '
$ws.Range('B50').Value = 'Stacks'
$ws.Range('C50').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.82'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.64%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '26.36'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.27%  '
